$wb = $excel.ActiveWorkbook

# --- Sheet 1: "Add Devices Loop A" ---
$ws1 = $wb.Worksheets.Item("Add Devices Loop A")

# Update the User Story text (shared string content change)
$ws1.Range("B3").Value = "NGC -1287 and NGC-491/TC-126"

# Update numeric value change
$ws1.Range("F8").Value = 0.34

# Update view: select B4 (this also resets the previous B1 top-left scroll position)
$ws1.Activate()
$ws1.Range("B4").Select() | Out-Null

# --- Sheet 2: "Delete Devices Loop A" ---
$ws2 = $wb.Worksheets.Item("Delete Devices Loop A")

# Update numeric value change
$ws2.Range("D7").Value = 0.54

# Update view: select C2:D2 with active cell C2
$ws2.Activate()
$ws2.Range("C2:D2").Select() | Out-Null

# Re-activate sheet 1 so it remains the selected/visible tab
$ws1.Activate()
